$d = $word.ActiveDocument

$replacements = @(
    @("303÷3=", "499÷5="),
    @("748÷5=", "585÷2="),
    @("223÷2=", "340÷4="),
    @("637÷7=", "944÷6="),
    @("690÷9=", "429÷8="),
    @("929÷5=", "915÷8="),
    @("466÷8=", "807÷9="),
    @("517÷3=", "462÷9="),
    @("116÷9=", "871÷4="),
    @("399÷7=", "729÷3="),
    @("985÷4=", "299÷3="),
    @("468÷6=", "125÷3="),
    @("213÷9=", "767÷7="),
    @("795÷2=", "778÷4="),
    @("433÷3=", "492÷4="),
    @("711÷4=", "695÷7="),
    @("488÷5=", "538÷2="),
    @("144÷5=", "571÷3="),
    @("537÷6=", "188÷7="),
    @("509÷3=", "793÷9="),
    @("900÷6=", "264÷3="),
    @("688÷3=", "239÷7="),
    @("495÷5=", "417÷4="),
    @("465÷3=", "290÷5="),
    @("303÷2=", "433÷6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
